# Update loading_percent values for the 380 kV case (Case_5_136).
# Data rows 2-25 correspond to line indices 0-23 (column A); the updated
# columns are B,C,D,E,F,H,I,K,M,N,O (columns G, J and L remain 0 and are
# left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$columns  = @(2,3,4,5,6,8,9,11,13,14,15)

$newValues = @(
    @(8.220758201098169,4.934674538863996,4.881175997749221,12.67062366166256,23.90263805532792,7.344005520526261,20.26921622773892,8.147799318923429,13.46619406058953,18.33807785221151,21.41567170606818),
    @(7.928069619367016,4.738234524528077,4.820547725536432,12.45736486036667,23.90998162060184,7.344005520526261,20.34139276162715,7.951884181708481,13.29288771035228,18.39639283048617,21.46902490471118),
    @(7.744015815610854,4.612328853306956,4.782338047141677,12.32896841116273,23.92074991205238,7.344005520526261,20.38948115228706,7.830426438291865,13.18861510946568,18.43385432712225,21.50634686456365),
    @(7.668040314494499,4.55973332520012,4.766527833758613,12.27736165416737,23.92671142334075,7.344005520526261,20.41002495245291,7.780711741322239,13.14671187391404,18.44953781459838,21.5227009766885),
    @(7.655369648301952,4.550923444353474,4.763888331842348,12.2688378178925,23.92779631768748,7.344005520526261,20.41349342447841,7.772445587380145,13.13979089846949,18.45216730573743,21.52548565680573),
    @(7.742994952898908,4.611624686104232,4.782125784066315,12.32826942670869,23.92082394218289,7.344005520526261,20.38975437859705,7.829756758884692,13.18804753758577,18.43406414753481,21.50656278846699),
    @(8.120808665561281,4.86806380884472,4.860480196921204,12.59660952864928,23.9038709046269,7.344005520526261,20.29331926060139,8.080532563724455,13.40602517712248,18.35784193395293,21.43311955355918),
    @(8.822558160893996,5.32737697964131,5.005993962240006,13.1396533508682,23.92027909411684,7.344005520526261,20.13418302714154,8.560009463997263,13.84815894264731,18.22145116971604,21.32539691432499),
    @(9.30875641816662,5.636506771013741,5.107512961311476,13.54438947012201,23.96253064923669,7.344005520526261,20.03560008479004,8.900969501184566,14.17885303274382,18.12913688263252,21.2685023827634),
    @(9.522602549761711,5.770699597201629,5.152440654569376,13.72883384716153,23.98827233850565,7.344005520526261,19.99474654074803,9.052899603146656,14.32991228314893,18.08883679835117,21.24747247171424),
    @(9.602459485728085,5.820570788022823,5.169266447463224,13.79865038772125,23.99895294670428,7.344005520526261,19.97985158671253,9.109920524694054,14.38715058152368,18.07381852327941,21.24020805139578),
    @(9.585311760542185,5.809872433777523,5.165651153475741,13.78361659136023,23.99661129065876,7.344005520526261,19.98303387008165,9.097663640595318,14.37482256752967,18.07704221041464,21.24174145833726),
    @(9.529195273929741,5.774821566049816,5.153828716786197,13.73457869164323,23.98913237620671,7.344005520526261,19.99350958470745,9.057601305006182,14.33462090967065,18.08759638232139,21.24686080382911),
    @(9.494674315882415,5.75322833265585,5.14656253024006,13.70453564594264,23.98467263681506,7.344005520526261,20.00000123439892,9.032993762969799,14.30999925345002,18.09409266380938,21.25008764069554),
    @(9.29462758615869,5.627605760554173,5.104550991384119,13.53233566204101,23.96097919354832,7.344005520526261,20.03835040775585,8.890971813416053,14.16898885208931,18.13180458217717,21.26997449764673),
    @(9.169977243428573,5.548877887368683,5.078452030784126,13.42672583092447,23.9481110863769,7.344005520526261,20.06289984857955,8.802992225524189,14.08260451854616,18.15537273149119,21.28341820673368),
    @(9.097595428052143,5.502991265674559,5.063322987879807,13.3660170839505,23.94132386531494,7.344005520526261,20.0773957377154,8.752092934751678,14.03298054403685,18.16908803538276,21.29160726349669),
    @(9.072972548324804,5.487351689384895,5.058180573839182,13.34547054718782,23.93913144773719,7.344005520526261,20.08236828355327,8.734810267819636,14.01619104507913,18.1737592350967,21.2944583082689),
    @(9.183318029396208,5.557321313964807,5.081242534019253,13.43796508805359,23.94941738623828,7.344005520526261,20.06024762563135,8.812388826863481,14.0917942583305,18.15284735940281,21.28193983399479),
    @(9.545708997115664,5.785142637366151,5.157306392372447,13.74898367551843,23.99130384374588,7.344005520526261,19.99041698698755,9.069382890352621,14.3464285771307,18.08448979393588,21.2453381425967),
    @(9.775987070629341,5.928523438850548,5.205923443705226,13.95205287832074,24.02411303692286,7.344005520526261,19.94813307071268,9.264388958834932,14.5130290358799,18.04122708891403,21.22549252207374),
    @(9.653704424264108,5.852508654515669,5.180078061282051,13.84371353768598,24.00610685233274,7.344005520526261,19.97039343474963,9.147701070672364,14.42411222347823,18.06418829578033,21.23571114922773),
    @(9.177288893293584,5.553505985779399,5.079981334452832,13.43288378996894,23.94882490502092,7.344005520526261,20.0614455055278,8.808141611681068,14.08763944967552,18.15398856399847,21.28260677325934),
    @(8.637517193199709,5.207984684065021,4.967545162503132,12.99141156243946,23.91052810313497,7.344005520526261,20.17401833771684,8.432024242600129,13.72729739277819,18.25695665149566,21.35063941834082)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowValues = $newValues[$i]
    $row = $startRow + $i
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $col = $columns[$j]
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
